# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right before the "总计" (Total) sheet,
#    populated with the same fund-holdings layout used by the other
#    quarterly sheets.
# 2. Insert a new top data row into "总计" summarising the new quarter,
#    pushing the existing rows down.
#
# Runtime quirks this script works around:
#   * A PowerShell variable holding a Worksheet reference is bound to that
#     sheet's *position* in the Worksheets collection, not to the sheet
#     itself. Calling Worksheets.Add(...) shifts every sheet at/after the
#     insertion point, so any worksheet variable captured earlier can end
#     up silently pointing at a different sheet afterwards. To stay safe,
#     Add() is only called ONCE in this whole script, and every worksheet
#     needed afterwards is (re-)fetched by name only after that call.
#   * Typing a digit-only string (e.g. "001479" or "9.45") into a
#     General-formatted cell gets auto-coerced to a number, which would
#     strip the leading zero / trailing zero the source data relies on.
#     Pre-formatting a scratch cell as Text, copying it, and
#     "Paste Special -> Values" into the destination carries the literal
#     string over without stamping the destination with a text number
#     format (matching the source file, where those cells carry no
#     explicit style).

$excel.DisplayAlerts = $false
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" worksheet, inserted before "总计"
# ---------------------------------------------------------------------------
$totalSheetBefore = $wb.Worksheets.Item("总计")
$newSheet = $wb.Worksheets.Add($totalSheetBefore)
$newSheet.Name = "2022-Q1"

# Re-fetch every sheet handle by name now that the one and only Add() call
# has happened.
$totalSheet = $wb.Worksheets.Item("总计")
$fmtSrc = $wb.Worksheets.Item("2021-Q4")

# Match the outline + page-setup conventions used by every other sheet in
# this workbook (the default "fresh sheet" ones differ slightly).
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Pull the header-row / index-column formatting (bold, centered, boxed)
# from an existing quarterly sheet so the new one matches visually.
$fmtSrc.Range("B1:H1").Copy()
$newSheet.Range("B1").PasteSpecial(-4122)   # xlPasteFormats
$fmtSrc.Range("A2:A4").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# Scratch cell (well outside the sheet's real data) used to smuggle
# digit-looking strings in as literal text - see note above.
$scratch = $newSheet.Range("A100")
$scratch.NumberFormat = "@"

function Set-Text($range, [string]$text) {
    $scratch.Value = $text
    $scratch.Copy()
    $range.PasteSpecial(-4163)   # xlPasteValues
}

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Data rows
$newSheet.Range("A2").Value = 0
Set-Text $newSheet.Range("B2") "001479"
$newSheet.Range("C2").Value = "中邮风格轮动灵活配置混合"
Set-Text $newSheet.Range("D2") "9.45"
Set-Text $newSheet.Range("E2") "62.17"
Set-Text $newSheet.Range("F2") "3.19"
Set-Text $newSheet.Range("G2") "0.3015"
$newSheet.Range("H2").Value = 7

$newSheet.Range("A3").Value = 1
Set-Text $newSheet.Range("B3") "005335"
$newSheet.Range("C3").Value = "浙商全景消费混合"
Set-Text $newSheet.Range("D3") "2.30"
Set-Text $newSheet.Range("E3") "93.36"
Set-Text $newSheet.Range("F3") "5.92"
Set-Text $newSheet.Range("G3") "0.1362"
$newSheet.Range("H3").Value = 9

$newSheet.Range("A4").Value = 2
Set-Text $newSheet.Range("B4") "003981"
$newSheet.Range("C4").Value = "中银证券瑞益灵活配置混合C"
Set-Text $newSheet.Range("D4") "0.21"
Set-Text $newSheet.Range("E4") "89.21"
Set-Text $newSheet.Range("F4") "3.57"
Set-Text $newSheet.Range("G4") "0.0075"
$newSheet.Range("H4").Value = 7

# ---------------------------------------------------------------------------
# 2. Add the 2022-Q1 summary row at the top of "总计"
# ---------------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Keep the bold/boxed index-column styling consistent with the rows below.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
Set-Text $totalSheet.Range("B2") "2022-Q1"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.45

# Only NOW drop the scratch row (every Set-Text call is done, so reusing the
# same scratch cell repeatedly above never leaves residue behind on
# "2022-Q1" - deleting it any earlier would have been re-populated by the
# later Set-Text calls above).
$newSheet.Rows.Item(100).Delete()

# ---------------------------------------------------------------------------
# Restore the originally-active sheet/selection (adding the sheet moved it)
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
